# Update the "No Register" value from UP032301000071 to UP032303000171
# This affects:
#  - Cell N2 which stores the register number on its own
#  - Cell F2 which stores a multi-line summary whose last line embeds the
#    same register number
# Also updates the active selection on the sheet from F2 to G2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldRegister = "UP032301000071"
$newRegister = "UP032303000171"

# Update the standalone register-number cell (N2)
$ws.Range("N2").Value2 = $newRegister

# Update the multi-line preparation text cell (F2), replacing the register
# number embedded on the last line while keeping the rest identical.
$f2 = $ws.Range("F2").Value2
$f2 = $f2.Replace($oldRegister, $newRegister)
$ws.Range("F2").Value2 = $f2

# Update the active selection to G2
$ws.Range("G2").Select()
